# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# --- Sheet "Home win": drop oldest match, shift remaining rows up ---
$ws1 = $wb.Worksheets.Item("Home win")
$ws1.Range("A2").Value = "23-01-2025 17:45"
$ws1.Range("B2").Value = "WORLD"
$ws1.Range("C2").Value = "UEFA EUROPA LEAGUE"
$ws1.Range("D2").Value = "AZ Alkmaar - AS Roma"
$ws1.Range("E2").Value = 70
$ws1.Range("F2").Value = 3

$ws1.Range("A3").Value = "23-01-2025 17:00"
$ws1.Range("B3").Value = "PORTUGAL"
$ws1.Range("C3").Value = "LIGA REVELAÇÃO U23"
$ws1.Range("D3").Value = "Sporting CP U23 - Benfica U23"
$ws1.Range("E3").Value = 73.3
$ws1.Range("F3").Value = 2.5

$ws1.Rows.Item(4).Delete()

# --- Sheet "Draw": drop three oldest matches, keep the newest as row 2 ---
$ws2 = $wb.Worksheets.Item("Draw")
$ws2.Range("A2").Value = "23-01-2025 23:00"
$ws2.Range("B2").Value = "BRAZIL"
$ws2.Range("C2").Value = "SERGIPANO"
$ws2.Range("D2").Value = "Barra SE - Lagarto"
$ws2.Range("E2").Value = 60
$ws2.Range("F2").Value = 5.25

$ws2.Rows.Item(5).Delete()
$ws2.Rows.Item(4).Delete()
$ws2.Rows.Item(3).Delete()

# --- Sheet "Over_Under": drop oldest match, shift remaining rows up ---
$ws4 = $wb.Worksheets.Item("Over_Under")
$ws4.Range("A2").Value = "23-01-2025 17:45"
$ws4.Range("B2").Value = "WORLD"
$ws4.Range("C2").Value = "UEFA EUROPA LEAGUE"
$ws4.Range("D2").Value = "Bodo/Glimt - Maccabi Tel Aviv"
$ws4.Range("E2").Value = 93.3
$ws4.Range("F2").Value = 1.62
$ws4.Range("G2").Value = 66.7
$ws4.Range("H2").Value = 2.5

$ws4.Range("A3").Value = "23-01-2025 17:45"
$ws4.Range("B3").Value = "WORLD"
$ws4.Range("C3").Value = "UEFA EUROPA LEAGUE"
$ws4.Range("D3").Value = "Malmo FF - Twente"
$ws4.Range("E3").Value = 73.3
$ws4.Range("F3").Value = 1.73
$ws4.Range("G3").Value = 60
$ws4.Range("H3").Value = 2.75

$ws4.Range("A4").Value = "23-01-2025 20:00"
$ws4.Range("B4").Value = "WORLD"
$ws4.Range("C4").Value = "UEFA EUROPA LEAGUE"
$ws4.Range("D4").Value = "Manchester United - Rangers"
$ws4.Range("E4").Value = 80
$ws4.Range("F4").Value = 1.62
$ws4.Range("G4").Value = 60
$ws4.Range("H4").Value = 2.5

$ws4.Rows.Item(5).Delete()
